$d = $word.ActiveDocument

# --- Helper: the flat-OPC wrapper required by Range.InsertXML ---
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# -----------------------------------------------------------------
# 1) "Git init" paragraph: drop the spell-check proofErr wrapping so
#    the whole line becomes a single run "Git init".
# -----------------------------------------------------------------
$cr = [char]13
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -eq ("Git init" + $cr)) {
        $target = $cand
        break
    }
}
if ($null -eq $target) {
    throw "Could not locate paragraph with text 'Git init'"
}

$rng = $target.Range
$rng.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
$rng.Delete()
$rng.InsertXML($pkgHeader + '<w:p><w:r><w:t>Git init</w:t></w:r></w:p>' + $pkgFooter)

# -----------------------------------------------------------------
# 2) Insert the new block of paragraphs right after the blank
#    paragraph that follows "Git init" (and before the final,
#    underline-formatted paragraph).
# -----------------------------------------------------------------
$afterGitInit = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq ("Git init" + $cr)) {
        $afterGitInit = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($null -eq $afterGitInit) {
    throw "Could not locate the blank paragraph following 'Git init'"
}

$insertionPoint = $d.Range($afterGitInit.Range.End, $afterGitInit.Range.End)

$newBody = '<w:p><w:r><w:t xml:space="preserve">Adicionar </w:t></w:r><w:r><w:t xml:space="preserve">um </w:t></w:r><w:r><w:t>arquivo no reposit&#243;rio:</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t>Git add &#8220;nome do arquivo&#8221;</w:t></w:r></w:p>' + `
           '<w:p/>' + `
           '<w:p><w:r><w:t>Adicionar v&#225;rios arquivos no reposit&#243;rio:</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t>Git add .</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t xml:space="preserve">Adicionar coment&#225;rio </w:t></w:r><w:r><w:t xml:space="preserve">em um </w:t></w:r><w:r><w:t>arquivo adicionado ao reposit&#243;rio:</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t>Git  commit</w:t></w:r><w:r><w:t xml:space="preserve"> &#8220;nome do arquivo&#8221;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t>m &#8220;</w:t></w:r><w:r><w:t>Estou enviando somente nome do arquivo</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r></w:p>' + `
           '<w:p/>' + `
           '<w:p><w:r><w:t>Adicionar coment&#225;rio &#250;nico em todos arquivos que ser&#227;o adicionados ao reposit&#243;rio:</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t xml:space="preserve">Git commit </w:t></w:r><w:r><w:t xml:space="preserve">-a </w:t></w:r><w:r><w:t>-m &#8220;Estou enviando todos os arquivos&#8221;</w:t></w:r></w:p>' + `
           '<w:p/>' + `
           '<w:p><w:r><w:t>Git branch -m main</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t>Git branch -m master</w:t></w:r></w:p>' + `
           '<w:p/>' + `
           '<w:p><w:r><w:t xml:space="preserve">Git Push -u origin master  </w:t></w:r></w:p>' + `
           '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>'

$insertionPoint.InsertXML($pkgHeader + $newBody + $pkgFooter)

# -----------------------------------------------------------------
# 3) Normal style: single line spacing w:line="256" w:lineRule="auto"
# -----------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle
$normal.ParagraphFormat.LineSpacing = 12.8    # -> w:line="256" w:lineRule="auto"

Write-Host "Edit complete"
